# Update "想去人数" (interested-count) figures in the 广州-漫展信息 workbook.
# This mirrors a refreshed scrape: several F-column counts increased slightly.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 156
$ws1.Range("F6").Value  = 713
$ws1.Range("F7").Value  = 1268
$ws1.Range("F9").Value  = 880
$ws1.Range("F11").Value = 275
$ws1.Range("F15").Value = 1063
$ws1.Range("F16").Value = 12241
$ws1.Range("F17").Value = 667
$ws1.Range("F20").Value = 48
$ws1.Range("F22").Value = 298
$ws1.Range("F23").Value = 1813
$ws1.Range("F26").Value = 501
$ws1.Range("F27").Value = 196
$ws1.Range("F28").Value = 113
$ws1.Range("F31").Value = 279
$ws1.Range("F32").Value = 94
$ws1.Range("F33").Value = 106
$ws1.Range("F35").Value = 189
$ws1.Range("F37").Value = 1215
$ws1.Range("F38").Value = 50

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 161
$ws2.Range("F9").Value  = 257
$ws2.Range("F11").Value = 102
$ws2.Range("F20").Value = 5

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1268
$ws4.Range("F7").Value  = 880
$ws4.Range("F9").Value  = 275
$ws4.Range("F14").Value = 1063
$ws4.Range("F15").Value = 12241
$ws4.Range("F16").Value = 257
$ws4.Range("F17").Value = 667
$ws4.Range("F20").Value = 48
$ws4.Range("F21").Value = 298
$ws4.Range("F22").Value = 1813
$ws4.Range("F24").Value = 501
$ws4.Range("F25").Value = 196
$ws4.Range("F26").Value = 102
$ws4.Range("F27").Value = 102
$ws4.Range("F29").Value = 113
$ws4.Range("F36").Value = 279
$ws4.Range("F37").Value = 94
$ws4.Range("F38").Value = 106
$ws4.Range("F41").Value = 189
$ws4.Range("F45").Value = 1215
$ws4.Range("F46").Value = 5
